$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Worksheet")

# Row 3: identifier changes from numeric 533 to the text "6464" (kept as
# text, not a number), and the "Save as..." label becomes "Pug".
# Force text storage for the numeric-looking value, then restore the
# original (General) cell formatting so the cell's style is unaffected.
$ws.Range("A3").NumberFormat = "@"
$ws.Range("A3").Value = "6464"
$ws.Range("A2").Copy()
$ws.Range("A3").PasteSpecial(-4122)
$ws.Range("B3").Value = "Pug"

# Row 11: add the missing German translation for "Password recover".
$ws.Range("C11").Value = "Passwort-Wiederherstellung"

# Row 12: add the missing German translation for
# "Your password reseted successully!".
$ws.Range("C12").Value = "Ihr Passwort wurde erfolgreich zurückgesetzt!"
